$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("http://dbpedia.org/resource/George_Murray_(bishop_of_Rochester)", "http://dbpedia.org/resource/Chester_Square", "http://dbpedia.org/resource/Lord_George_Murray_(bishop)", "http://dbpedia.org/resource/1860"),
    @("http://dbpedia.org/resource/Geoffrey_(archbishop_of_York)", "http://dbpedia.org/resource/Normandy", "http://dbpedia.org/resource/Henry_II_of_England", "http://dbpedia.org/resource/12-12-12"),
    @("http://dbpedia.org/resource/Ferdinand_III_of_Castile", "http://dbpedia.org/resource/Crown_of_Castile", "http://dbpedia.org/resource/Alfonso_IX_of_León", "http://dbpedia.org/resource/1252"),
    @("http://dbpedia.org/resource/Erik_Benzelius_the_younger", "http://dbpedia.org/resource/Linköping", "http://dbpedia.org/resource/Erik_Benzelius_the_Elder", "http://dbpedia.org/resource/1743"),
    @("http://dbpedia.org/resource/Edward_the_Confessor", "http://dbpedia.org/resource/London", "http://dbpedia.org/resource/Æthelred_the_Unready", "http://dbpedia.org/resource/1066"),
    @("http://dbpedia.org/resource/Edward_William_Grinfield", "http://dbpedia.org/resource/Brighton", "http://dbpedia.org/resource/Thomas_Grinfield", "http://dbpedia.org/resource/1864"),
    @("http://dbpedia.org/resource/Edward_Francis_Wilson", "http://dbpedia.org/resource/Saltspring_Island", "http://dbpedia.org/resource/Daniel_Wilson_(bishop)", "http://dbpedia.org/resource/1915"),
    @("http://dbpedia.org/resource/Donald_Foster_Hudson", "http://dbpedia.org/resource/England", "http://dbpedia.org/resource/Father", "http://dbpedia.org/resource/2003")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
